$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($row, $col, $oldText, $newText) {
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    # Exclude the trailing cell-mark character so we only replace the visible text
    $rng.MoveEnd(1, -1) | Out-Null
    if ($rng.Text -ne $oldText) {
        Write-Host "WARNING: Cell($row,$col) expected [$oldText] but found [$($rng.Text)]"
    }
    $rng.Text = $newText
}

Set-CellText 1 1 "59×45=" "73×83="
Set-CellText 1 2 "67×95=" "70×43="
Set-CellText 1 3 "96×80=" "40×52="
Set-CellText 1 4 "27×26=" "56×77="
Set-CellText 1 5 "99×38=" "87×93="
Set-CellText 5 1 "99×46=" "50×30="
Set-CellText 5 2 "62×82=" "76×59="
Set-CellText 5 3 "61×56=" "22×60="
Set-CellText 5 4 "55×79=" "49×23="
Set-CellText 5 5 "92×95=" "73×23="
Set-CellText 10 1 "94×89=" "52×55="
Set-CellText 10 2 "28×62=" "82×96="
Set-CellText 10 3 "20×88=" "15×53="
Set-CellText 10 4 "86×75=" "93×46="
Set-CellText 10 5 "60×99=" "25×26="
Set-CellText 15 1 "38×99=" "35×38="
Set-CellText 15 2 "18×83=" "71×11="
Set-CellText 15 3 "60×98=" "82×55="
Set-CellText 15 4 "16×97=" "91×45="
Set-CellText 15 5 "26×14=" "48×37="
Set-CellText 20 1 "18×30=" "26×44="
Set-CellText 20 2 "16×30=" "24×61="
Set-CellText 20 3 "73×23=" "92×88="
Set-CellText 20 4 "73×60=" "57×37="
Set-CellText 20 5 "27×41=" "29×35="
